$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix B35: it was mistakenly entered as text "3"; it should be the number 3 ---
$ws.Range("B35").Value = 3

# --- Append a new annotation row (row 36) ---
$ws.Range("A36").Value = "Sunsi Wu"

# B36's politeness_score of "2" needs to be stored as text (matching the same
# quirk B35 had), not a number. Build it via a TEXT() formula, then copy/paste
# the computed value so the literal lands as a text cell without leaving a
# formula behind or minting a new (text) number-format style on the cell.
$ws.Range("B36").Formula = "=TEXT(2,""0"")"
$ws.Range("B36").Copy()
$ws.Range("B36").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("C36").Value = "无"
$ws.Range("D36").Value = "ACK"
$ws.Range("E36").Value = "WRI"
$ws.Range("F36").Value = "d3d18e54-3fa0-4f4c-98ec-51cec5852681"
$ws.Range("G36").Value = "HknbyQbC-_annotated.xlsx"
$ws.Range("H36").Value = "- Added suggested references and updated section 2 to include more comprehensive analysis for related work."
